# Auto-applied numeric updates to the Leve profit tables (ALC, ARM, BSM, CRP, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 115.35714
$ws.Range("I9").Value = 65.375
$ws.Range("J9").Value = 182
$ws.Range("K9").Value = 65.375
$ws.Range("L9").Value = 182
$ws.Range("M9").Value = 103.625
$ws.Range("N9").Value = -520

$ws.Range("H15").Value = 1474.1976
$ws.Range("I15").Value = 1474.1976
$ws.Range("K15").Value = 4422.5928
$ws.Range("M15").Value = -4253.5928

$ws.Range("H40").Value = 1262.7805
$ws.Range("I40").Value = 1199.9166
$ws.Range("J40").Value = 1351.5294
$ws.Range("K40").Value = 1199.9166
$ws.Range("L40").Value = 1351.5294
$ws.Range("M40").Value = -1024.9166
$ws.Range("N40").Value = -1701.5294

$ws.Range("H64").Value = 3164.4827
$ws.Range("I64").Value = 3099.2593
$ws.Range("J64").Value = 4045
$ws.Range("K64").Value = 3099.2593
$ws.Range("L64").Value = 4045
$ws.Range("M64").Value = -2851.2593
$ws.Range("N64").Value = -4541

$ws.Range("H67").Value = 3164.4827
$ws.Range("I67").Value = 3099.2593
$ws.Range("J67").Value = 4045
$ws.Range("K67").Value = 3099.2593
$ws.Range("L67").Value = 4045
$ws.Range("M67").Value = -2241.2593
$ws.Range("N67").Value = -5761

$ws.Range("H76").Value = 5294433.5
$ws.Range("I76").Value = 6176172.5
$ws.Range("J76").Value = 4000
$ws.Range("K76").Value = 6176172.5
$ws.Range("L76").Value = 4000
$ws.Range("M76").Value = -6175857.5
$ws.Range("N76").Value = -4630

$ws.Range("H79").Value = 5294433.5
$ws.Range("I79").Value = 6176172.5
$ws.Range("J79").Value = 4000
$ws.Range("K79").Value = 6176172.5
$ws.Range("L79").Value = 4000
$ws.Range("M79").Value = -6175080.5
$ws.Range("N79").Value = -6184

$ws.Range("H88").Value = 6164075
$ws.Range("I88").Value = 1189.5555
$ws.Range("J88").Value = 10786239
$ws.Range("K88").Value = 1189.5555
$ws.Range("L88").Value = 10786239
$ws.Range("M88").Value = -783.5554999999999
$ws.Range("N88").Value = -10787051

$ws.Range("H91").Value = 6164075
$ws.Range("I91").Value = 1189.5555
$ws.Range("J91").Value = 10786239
$ws.Range("K91").Value = 1189.5555
$ws.Range("L91").Value = 10786239
$ws.Range("M91").Value = 214.4445000000001
$ws.Range("N91").Value = -10789047

$ws.Range("H100").Value = 25644064
$ws.Range("I100").Value = 37039092
$ws.Range("J100").Value = 5250
$ws.Range("K100").Value = 37039092
$ws.Range("L100").Value = 5250
$ws.Range("M100").Value = -37038551
$ws.Range("N100").Value = -6332

$ws.Range("H103").Value = 100000800
$ws.Range("I103").Value = 1000
$ws.Range("J103").Value = 125000750
$ws.Range("K103").Value = 3000
$ws.Range("L103").Value = 375002250
$ws.Range("M103").Value = -2414
$ws.Range("N103").Value = -375003422

$ws.Range("H106").Value = 2666.6667
$ws.Range("I106").Value = 2666.6667
$ws.Range("K106").Value = 2666.6667
$ws.Range("M106").Value = -2035.6667

$ws.Range("H107").Value = 3666.6667
$ws.Range("J107").Value = 3666.6667
$ws.Range("L107").Value = 3666.6667
$ws.Range("N107").Value = -7506.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2477.44
$ws.Range("I45").Value = 2173.1428
$ws.Range("J45").Value = 2864.7273
$ws.Range("K45").Value = 2173.1428
$ws.Range("L45").Value = 2864.7273
$ws.Range("M45").Value = -1796.1428
$ws.Range("N45").Value = -3618.7273

$ws.Range("H55").Value = 33132.555
$ws.Range("J55").Value = 33132.555
$ws.Range("L55").Value = 33132.555
$ws.Range("N55").Value = -33762.555

$ws.Range("H63").Value = 7666.3335
$ws.Range("I63").Value = 7499.5
$ws.Range("K63").Value = 7499.5
$ws.Range("M63").Value = -6813.5

$ws.Range("H66").Value = 7666.3335
$ws.Range("I66").Value = 7499.5
$ws.Range("K66").Value = 37497.5
$ws.Range("M66").Value = -34065.5

$ws.Range("H80").Value = 20666.934
$ws.Range("J80").Value = 20666.934
$ws.Range("L80").Value = 20666.934
$ws.Range("N80").Value = -22662.934

$ws.Range("H83").Value = 20666.934
$ws.Range("J83").Value = 20666.934
$ws.Range("L83").Value = 62000.802
$ws.Range("N83").Value = -71984.802

$ws.Range("H140").Value = 116000
$ws.Range("J140").Value = 116000
$ws.Range("L140").Value = 116000
$ws.Range("N140").Value = -126360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()

$ws.Range("H31").Value = 3254.3333
$ws.Range("J31").Value = 3254.3333
$ws.Range("L31").Value = 3254.3333
$ws.Range("N31").Value = -3758.3333

$ws.Range("H35").Value = 11090.909
$ws.Range("J35").Value = 11090.909
$ws.Range("L35").Value = 11090.909
$ws.Range("N35").Value = -11710.909

$ws.Range("H82").Value = 15024.6
$ws.Range("I82").Value = 6185.6665
$ws.Range("J82").Value = 28283
$ws.Range("K82").Value = 6185.6665
$ws.Range("L82").Value = 28283
$ws.Range("M82").Value = -5802.6665
$ws.Range("N82").Value = -29049

$ws.Range("H85").Value = 15024.6
$ws.Range("I85").Value = 6185.6665
$ws.Range("J85").Value = 28283
$ws.Range("K85").Value = 6185.6665
$ws.Range("L85").Value = 28283
$ws.Range("M85").Value = -4859.6665
$ws.Range("N85").Value = -30935

$ws.Range("H140").Value = 72666.336
$ws.Range("J140").Value = 72666.336
$ws.Range("L140").Value = 72666.336
$ws.Range("N140").Value = -83026.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 217.5
$ws.Range("I12").Value = 217.5
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 217.5
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -47.5
$ws.Range("N12").ClearContents()

$ws.Range("H41").Value = 5744.3335
$ws.Range("I41").Value = 3333
$ws.Range("J41").Value = 6950
$ws.Range("K41").Value = 3333
$ws.Range("L41").Value = 6950
$ws.Range("M41").Value = -2905
$ws.Range("N41").Value = -7806

$ws.Range("H50").Value = 8000
$ws.Range("J50").Value = 8000
$ws.Range("L50").Value = 8000
$ws.Range("N50").Value = -9250

$ws.Range("H59").Value = 11350.875
$ws.Range("I59").Value = 10000
$ws.Range("J59").Value = 11801.167
$ws.Range("K59").Value = 10000
$ws.Range("L59").Value = 11801.167
$ws.Range("M59").Value = -8855
$ws.Range("N59").Value = -14091.167

$ws.Range("H74").Value = 14000
$ws.Range("I74").Value = 14000
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 14000
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -13126
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 14000
$ws.Range("I77").Value = 14000
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 42000
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -37632
$ws.Range("N77").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H41").Value = 1670066.4

$ws.Range("H56").Value = 10017.286
$ws.Range("J56").Value = 11014
$ws.Range("L56").Value = 11014
$ws.Range("N56").Value = -12396

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 47714.5
$ws.Range("J46").Value = 47714.5
$ws.Range("L46").Value = 47714.5
$ws.Range("N46").Value = -48176.5

$ws.Range("H134").Value = 47714.5
$ws.Range("J134").Value = 47714.5
$ws.Range("L134").Value = 143143.5
$ws.Range("N134").Value = -148213.5

